$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark; we'll re-insert a fresh one
# in the correct position as part of the new content below.
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
}

$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.01</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Erros de sintaxe: erros de gramática e pontuação</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">              </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.01.01</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Recomendação: Verificar ortografia</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.02</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Erros de runtime: aparecem à medida que o programa é executado</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">              </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.02.01</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Recomendação: Rever o código, por exemplo, procurar divisões por zero</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.03</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>– Erros lógicos: erros básicos em algoritmos ou procedimentos, apenas detectados quando o resultado destes é anomalo.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">              </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.03.01</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Recomendação: Usar mecanismos de debug (passo a passo)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.04</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Erros de thre</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>ding:  Ocorrem quando se usam threads, exemplos podem ser deadlocks e concorrência no acesso.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">              </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.04.01</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Recomendação:  Utilizar mecanismos de defesa, por exemplo, semáforos.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.05</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Excepções: Não são erros fatais mas podem provocar problemas na execução do programa</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">               </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PARS.02.05.01</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> – Recomendação: Utilizar técnicas de controlo de excepções, por exemplo, “try catch”</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>'

$r.InsertXML($xml)

Write-Output "done"
